$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 100, shifting existing rows 100:219 down to 102:221
$ws.Rows.Item(100).Resize(2).Insert()

# --- New row 100 ---
$ws.Cells.Item(100,1).Value = 3
$ws.Cells.Item(100,2).Value = "Femacal de La Calera"
$ws.Cells.Item(100,3).Value = "Coquimbo"
$ws.Cells.Item(100,4).Value = 44483
$ws.Cells.Item(100,5).Value = 5
$ws.Cells.Item(100,6).Value = 100112040
$ws.Cells.Item(100,7).Value = "Cilantro"
$ws.Cells.Item(100,8).Value = "Sin especificar"
$ws.Cells.Item(100,9).Value = "Primera"
$ws.Cells.Item(100,10).Value = 180
$ws.Cells.Item(100,11).Value = 2500
$ws.Cells.Item(100,12).Value = 2500
$ws.Cells.Item(100,13).Value = 2500
$ws.Cells.Item(100,14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(100,15).Value = "Provincia de Quillota"
$ws.Cells.Item(100,16).Value = 833
$ws.Cells.Item(100,17).Value = 3
$ws.Cells.Item(100,18).Value = "Hortaliza"

# --- New row 101 ---
$ws.Cells.Item(101,1).Value = 3
$ws.Cells.Item(101,2).Value = "Femacal de La Calera"
$ws.Cells.Item(101,3).Value = "Coquimbo"
$ws.Cells.Item(101,4).Value = 44483
$ws.Cells.Item(101,5).Value = 5
$ws.Cells.Item(101,6).Value = 100112040
$ws.Cells.Item(101,7).Value = "Cilantro"
$ws.Cells.Item(101,8).Value = "Sin especificar"
$ws.Cells.Item(101,9).Value = "Segunda"
$ws.Cells.Item(101,10).Value = 180
$ws.Cells.Item(101,11).Value = 2300
$ws.Cells.Item(101,12).Value = 2300
$ws.Cells.Item(101,13).Value = 2300
$ws.Cells.Item(101,14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(101,15).Value = "Provincia de Quillota"
$ws.Cells.Item(101,16).Value = 767
$ws.Cells.Item(101,17).Value = 3
$ws.Cells.Item(101,18).Value = "Hortaliza"
